# Fix Mining data: corrected total movement values from chart re-read
# Revised total movement against y-axis gridlines (0,20,40,60,80), with
# ore/waste recalculated from corrected totals + strip ratios.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Mining" - full Ore/Waste/Total Movement/Strip Ratio/Ore%/Waste% table
# ---------------------------------------------------------------------------
$wsMining = $wb.Worksheets.Item("Mining")

# Row 6 - 2028
$wsMining.Range("B6").Value = 10.2
$wsMining.Range("C6").Value = 31.8
$wsMining.Range("D6").Value = 42
$wsMining.Range("F6").Value = 24.3
$wsMining.Range("G6").Value = 75.7

# Row 7 - 2030
$wsMining.Range("B7").Value = 9.4
$wsMining.Range("C7").Value = 40.6
$wsMining.Range("D7").Value = 50
$wsMining.Range("F7").Value = 18.8
$wsMining.Range("G7").Value = 81.2

# Row 8 - 2032
$wsMining.Range("B8").Value = 8.800000000000001
$wsMining.Range("C8").Value = 51.2
$wsMining.Range("D8").Value = 60

# Row 9 - 2034 (peak mining year)
$wsMining.Range("B9").Value = 6
$wsMining.Range("C9").Value = 69
$wsMining.Range("D9").Value = 75

# Row 10 - 2036
$wsMining.Range("B10").Value = 8.300000000000001
$wsMining.Range("C10").Value = 41.7
$wsMining.Range("D10").Value = 50
$wsMining.Range("F10").Value = 16.6
$wsMining.Range("G10").Value = 83.40000000000001

# Row 11 - 2038
$wsMining.Range("B11").Value = 10
$wsMining.Range("C11").Value = 45
$wsMining.Range("D11").Value = 55

# Row 12 - 2040
$wsMining.Range("B12").Value = 6.9
$wsMining.Range("C12").Value = 38.1
$wsMining.Range("D12").Value = 45
$wsMining.Range("F12").Value = 15.3
$wsMining.Range("G12").Value = 84.7

# Row 13 - 2042
$wsMining.Range("C13").Value = 24
$wsMining.Range("D13").Value = 31
$wsMining.Range("F13").Value = 22.6
$wsMining.Range("G13").Value = 77.40000000000001

# Row 14 - 2044
$wsMining.Range("B14").Value = 6.2
$wsMining.Range("C14").Value = 11.8
$wsMining.Range("D14").Value = 18
$wsMining.Range("F14").Value = 34.4
$wsMining.Range("G14").Value = 65.59999999999999

# Row 16 - 2048
$wsMining.Range("C16").Value = 3.6
$wsMining.Range("D16").Value = 8
$wsMining.Range("F16").Value = 55
$wsMining.Range("G16").Value = 45

# Row 17 - LoM TOTAL
$wsMining.Range("B17").Value = 92.2
$wsMining.Range("C17").Value = 395.8
$wsMining.Range("D17").Value = 488
$wsMining.Range("E17").Value = 4.3
$wsMining.Range("F17").Value = 18.9
$wsMining.Range("G17").Value = 81.09999999999999

# ---------------------------------------------------------------------------
# Sheet "Summary & Analysis" - Ore Mined / Waste Mined / Total Movement columns
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary & Analysis")

# Row 6 - 2028
$wsSummary.Range("B6").Value = 10.2
$wsSummary.Range("C6").Value = 31.8
$wsSummary.Range("D6").Value = 42

# Row 7 - 2030
$wsSummary.Range("B7").Value = 9.4
$wsSummary.Range("C7").Value = 40.6
$wsSummary.Range("D7").Value = 50

# Row 8 - 2032
$wsSummary.Range("B8").Value = 8.800000000000001
$wsSummary.Range("C8").Value = 51.2
$wsSummary.Range("D8").Value = 60

# Row 9 - 2034
$wsSummary.Range("B9").Value = 6
$wsSummary.Range("C9").Value = 69
$wsSummary.Range("D9").Value = 75

# Row 10 - 2036
$wsSummary.Range("B10").Value = 8.300000000000001
$wsSummary.Range("C10").Value = 41.7
$wsSummary.Range("D10").Value = 50

# Row 11 - 2038
$wsSummary.Range("B11").Value = 10
$wsSummary.Range("C11").Value = 45
$wsSummary.Range("D11").Value = 55

# Row 12 - 2040
$wsSummary.Range("B12").Value = 6.9
$wsSummary.Range("C12").Value = 38.1
$wsSummary.Range("D12").Value = 45

# Row 13 - 2042
$wsSummary.Range("C13").Value = 24
$wsSummary.Range("D13").Value = 31

# Row 14 - 2044
$wsSummary.Range("B14").Value = 6.2
$wsSummary.Range("C14").Value = 11.8
$wsSummary.Range("D14").Value = 18

# Row 16 - 2048
$wsSummary.Range("C16").Value = 3.6
$wsSummary.Range("D16").Value = 8

# ---------------------------------------------------------------------------
# Sheet "Key Insights" - narrative text updated to match corrected totals
# ---------------------------------------------------------------------------
$wsInsights = $wb.Worksheets.Item("Key Insights")

$wsInsights.Range("B5").Value = "2034 - Total movement ~75.0 Mt with highest strip ratio of 11.5. This is a massive waste stripping campaign."
$wsInsights.Range("B6").Value = "92.2 Mt across even-year snapshots"
$wsInsights.Range("B7").Value = "395.8 Mt across even-year snapshots"
$wsInsights.Range("B8").Value = "4.29 (waste:ore) - Very high waste burden, especially in 2034"
